$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.129.56'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.41%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.656.04'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.29%  '

$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.02%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5307'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.71%  '

$ws.Range("E7").Value = '  +0.12%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2617'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.15%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06331'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.97%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.47'
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07802'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.77%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.518'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.19%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.656.89'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.30%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.884.32'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.36%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5498'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.52%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8197'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.44%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.44'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.71%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.138.26'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.41%  '

$ws.Range("E19").Value = '  +0.07%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.613'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.08%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '191.05'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.24%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.09'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.95%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.020'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.82%  '

$ws.Range("E24").Value = '  +0.16%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.35'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.04%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1230'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.37%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.224'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.35%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.98'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.99%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.468'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.49%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05712'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.87%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.273'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.03%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.560'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.31%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.267'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.26%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.598'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.09%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.802'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.81%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9516'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.90%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.416'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.19%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5733'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.41%  '

$ws.Range("E39").Value = '  +0.76%  '

$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8529'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.18%  '

$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.799'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.93%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '104.46'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.87%  '

$ws.Range("E43").Value = '  +0.17%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.037.70'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.46%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.797.72'
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.76'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.33%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.004'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.13%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4348'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.26%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.855'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.74%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05154'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.09%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.440'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.18%  '
